$d = $word.ActiveDocument

# --- Find the paragraph holding the M2Doc "for" field (it begins with the field). ---
$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $targetParagraph = $p
    }
}

$insPos = $targetParagraph.Range.Start

# --- The field carries a "_GoBack" bookmark (placed right after "eClassifiers",
#     right before the trailing space) that must survive the rewrite, so drop it
#     for now and put an equivalent one back once the field is gone. ---
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$field = $targetParagraph.Range.Fields(1)
$field.Delete()

# --- Re-insert the field's instruction text as literal, visible text, split into
#     the same four chunks the field used to be made of. ---
$part1 = "{m"
$part2 = ":for v| self."
$part3 = "eClassifiers"
$part4 = "}"
$fullText = $part1 + $part2 + $part3 + $part4

$insertionPoint = $d.Range($insPos, $insPos)
$insertionPoint.InsertBefore($fullText)

$p1End = $insPos + $part1.Length
$p2End = $p1End + $part2.Length
$p3End = $p2End + $part3.Length
$p4End = $p3End + $part4.Length

# --- Force the same run layout a real edit would leave behind: split the run
#     after "{m" and after ":for v| self." using a scratch bookmark that is
#     immediately removed again (adding/removing a bookmark breaks a run in two
#     without touching the surrounding text or its formatting). ---
$d.Bookmarks.Add("ScratchSplit1", $d.Range($p1End, $p1End)) | Out-Null
$d.Bookmarks("ScratchSplit1").Delete()

$d.Bookmarks.Add("ScratchSplit2", $d.Range($p2End, $p2End)) | Out-Null
$d.Bookmarks("ScratchSplit2").Delete()

# --- Put the _GoBack bookmark back where it originally sat: right after
#     "eClassifiers" and right before the closing "}". ---
$d.Bookmarks.Add("_GoBack", $d.Range($p3End, $p3End)) | Out-Null

# --- Keep "}" from re-merging into the following "A paragraph" run. ---
$d.Bookmarks.Add("ScratchSplit3", $d.Range($p4End, $p4End)) | Out-Null
$d.Bookmarks("ScratchSplit3").Delete()

Write-Output "Paragraph now reads: $($targetParagraph.Range.Text)"
